# Auto-update draw results: append the 2025-11-20 "Pick 3" draw as a new
# row right after the current last row of data (diff: row 65 added,
# dimension/used-range grows from A1:E64 to A1:E65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row so the new record lands directly below it
# (mirrors the bot's "append one row per day" behaviour).
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# All five columns in this sheet are plain text (Date/Game/Phase/Result/
# InsertedAt are stored as strings even when they look numeric, e.g.
# Phase "251120" or the date "2025-11-20"). Prefix the numeric-looking
# values with a leading apostrophe so Excel keeps them as text instead of
# auto-converting to a date serial / number.
$ws.Cells.Item($newRow, 1).Value = "'2025-11-20"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "'251120"
$ws.Cells.Item($newRow, 4).Value = "6-0-8"
$ws.Cells.Item($newRow, 5).Value = "'2025-11-20T21:38:40.871+04:00"

# Re-apply the plain (unstyled) formatting from the previous data row so
# the new row doesn't pick up the quote-prefix formatting mark and instead
# matches the look of every other row in the table.
$ws.Range("A" + $lastRow + ":E" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":E" + $newRow).PasteSpecial(-4122)
